$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert a new row at position 16 (old row 16 "pret.jsp / prolongerPret_ajax" shifts to row 17)
$ws.Rows.Item(16).Insert()

# 2. Bring over formatting for the new row 16 from row 15 (same visual style family: s5,s11,s3,s5,s9,s9,s8,s5)
$ws.Range("A15:H15").Copy()
$ws.Range("A16:H16").PasteSpecial(-4122)

# 3. A16 is the bottom of a new vertical merge (A15:A16), so it must use the same "bottom of merge" style as A10 (bottom of A4:A10)
$ws.Range("A10").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4. Merge A15:A16 (new "recherche.jsp" block spanning the lancerRecherche_ajax + voirDispo_ajax rows)
$ws.Range("A15:A16").Merge()

# 5. Row heights: row 16 is a new 45pt row, row 15 stays 60pt (untouched), row 17 stays 30pt (untouched, it just moved)
$ws.Rows.Item(16).RowHeight = 45

# 6. Fill in the new row 16 content (voirDispo_ajax) and the new D15 cell, in an order that matches
#    the append order of the new shared strings introduced by this change.
$ws.Cells.Item(16, 2).Value = "voirDispo_ajax"
$ws.Cells.Item(16, 3).Value = "Voir la dispo d'un livre dans le différente bibliothèque de la ville "
$ws.Cells.Item(15, 4).Value = "livre.RechercherAction / actionAjax"
$ws.Cells.Item(16, 6).Value = "(JSON) "
$ws.Cells.Item(16, 4).Value = "livre.DispoAction / actionAjax"
$ws.Cells.Item(17, 5).Value = "int id"

# E16 keeps the "int isbn" wording that used to live on the old row 16
$ws.Cells.Item(16, 5).Value = "int isbn"

# G16 keeps the same "RAS : action AJAX" wording as G15 (PasteSpecial only carried the formatting, not the value)
$ws.Cells.Item(16, 7).Value = "RAS : action AJAX"

# H16 / A16 stay blank, as they were for the row-15 paste source (D15's old blank cell, H15's blank cell).

# 7. Extend the AutoFilter range and the hidden _FilterDatabase defined name to the new last row (17)
$ws.AutoFilterMode = $false
$ws.Range("A1:H17").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name() -eq "Tabelle1!_FilterDatabase") {
        $n.RefersTo = "=Tabelle1!`$A`$1:`$H`$17"
    }
}

# 8. Update the view: scrolled so row 13 is at the top, with E18 selected (just below the edited rows)
$ws.Range("E18").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
